$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.361.08"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -0.04%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.913.07"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +0.42%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.722"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +9.49%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'254.74"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +4.24%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.25%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'40.83"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -2.00%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  +2.45%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'52.28"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -0.90%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.0754"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +5.55%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.0990"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -0.72%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'2.190.46"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.47%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'12.66"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +5.11%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.722"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +3.47%  "
$ws.Range("E15").ClearFormats()
$ws.Range("B16").Value = "'Polkadot"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C16").ClearFormats()
$ws.Range("D16").Value = "'4.94"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +1.57%  "
$ws.Range("E16").ClearFormats()
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("B17").ClearFormats()
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C17").ClearFormats()
$ws.Range("D17").Value = "'1.923.28"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +0.77%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'35.354.23"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -0.11%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'74.68"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +4.12%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'0.0₃0856"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +4.29%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'244.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +1.73%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'13.05"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +4.21%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'5.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +5.70%  "
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'  +0.16%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'2.48"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +4.10%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  +3.26%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'166.76"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -2.21%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'8.67"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +2.72%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'18.79"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +2.07%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  +5.22%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'4.130.38"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +19.50%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'4.37"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +5.35%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'2.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +14.45%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'1.64"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +21.06%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  +3.26%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'4.24"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +2.64%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  +0.10%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.921"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -1.50%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'2.04"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +0.49%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'17.37"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +6.15%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.0220"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +4.62%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'97.39"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +7.90%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  +1.69%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.0652"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +0.35%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'1.342.49"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -0.18%  "
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'  +1.72%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'2.43"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +0.78%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'6.78"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +2.96%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'  -0.74%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'45.32"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -5.52%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'11.76"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +8.03%  "
$ws.Range("E51").ClearFormats()
